$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ln_pgdp)
$ws.Cells.Item(2, 3).Value = -5.062293621080931
$ws.Cells.Item(2, 4).Value = 113.1898977236058
$ws.Cells.Item(2, 7).Value = -1.384053665424951
$ws.Cells.Item(2, 8).Value = 0.1664453888926849

# Row 3 (ln_pop_density)
$ws.Cells.Item(3, 3).Value = -1.676422489111038
$ws.Cells.Item(3, 4).Value = 104.6579459240546
$ws.Cells.Item(3, 7).Value = -0.4583413892060091
$ws.Cells.Item(3, 8).Value = 0.6467405249150434

# Row 4 (tertiary_share)
$ws.Cells.Item(4, 3).Value = 16.21594552480293
$ws.Cells.Item(4, 4).Value = 64.12677777939534
$ws.Cells.Item(4, 7).Value = 4.433511866730184
$ws.Cells.Item(4, 8).Value = 0.00000960515358271766

# Row 5 (tertiary_share_sq)
$ws.Cells.Item(5, 3).Value = 15.99776471493287
$ws.Cells.Item(5, 4).Value = 64.92584297378158
$ws.Cells.Item(5, 7).Value = 4.373860259725699
$ws.Cells.Item(5, 8).Value = 0.00001262426739122109

# Row 6 (ln_fdi)
$ws.Cells.Item(6, 3).Value = -6.273725120801471
$ws.Cells.Item(6, 4).Value = 119.2074657667888
$ws.Cells.Item(6, 7).Value = -1.715264443207027
$ws.Cells.Item(6, 8).Value = 0.08640043158271631

# Row 7 (ln_road_area)
$ws.Cells.Item(7, 3).Value = 0.3968688470246407
$ws.Cells.Item(7, 4).Value = -101.7303809120134
$ws.Cells.Item(7, 7).Value = 0.1085057137203632
$ws.Cells.Item(7, 8).Value = 0.9136018339016092
